# HRMS DATA ADDED SUCCESSFULLY
# Adds a new "LeaveConfiguration" sheet (business-unit -> weekend-day config)
# after the existing "BranchMaster" sheet, populates it, and makes it the
# active sheet/tab - mirroring how the other master-data sheets in this
# workbook are laid out.

$wb = $excel.ActiveWorkbook

# The sheet that used to be last/active (BranchMaster) - new sheet goes right after it.
$branchMaster = $wb.Worksheets.Item("BranchMaster")

$newSheet = $wb.Worksheets.Add($null, $branchMaster)
$newSheet.Name = "LeaveConfiguration"

# Header row
$newSheet.Range("A1").Value = "Bunit"
$newSheet.Range("B1").Value = "WEEK"
$newSheet.Range("C1").Value = "WeekEND"

# Data row
$newSheet.Range("A2").Value = "BU1-Test"
$newSheet.Range("B2").Value = "Saturday"
$newSheet.Range("C2").Value = "Sunday"

# Match the look of the workbook's other data sheets: center-aligned cells.
# Pull the formatting from an existing populated sheet (LacationMaster!A1)
# instead of poking HorizontalAlignment/VerticalAlignment directly so we
# reuse the existing shared cell style instead of minting a new one.
$styleSource = $wb.Worksheets.Item("LacationMaster").Range("A1")
$styleSource.Copy()
$newSheet.Range("A1:C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths close to the sibling master-data sheets.
$newSheet.Columns.Item(1).ColumnWidth = 17.7109375
$newSheet.Columns.Item(2).ColumnWidth = 23.42578125
$newSheet.Columns.Item(3).ColumnWidth = 21

# Leave the same cell selected/active as the source workbook.
$newSheet.Range("C11").Select()

# Make the new sheet the active tab (was BranchMaster before).
$newSheet.Activate()
